# Daily attendance processing - 2025-10-27 04:53:42
#
# Normalizes the "Recorded By" column (G) on the active sheet: whenever the
# cell holds a comma-separated list of recorders whose last entry is
# "System", the list order is reversed so that "System" is reported first.
# Rows whose "Recorded By" value is a single name, or which already starts
# with "System" in some other arrangement, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# "Recorded By" lives in column G (7).
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -notmatch ",") { continue }

    $parts = $value -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    if ($trimmed.Length -lt 2) { continue }
    if ($trimmed[$trimmed.Length - 1] -ne "System") { continue }

    $reversed = @()
    for ($i = $trimmed.Length - 1; $i -ge 0; $i--) { $reversed += $trimmed[$i] }

    $newValue = [string]::Join(", ", $reversed)
    $cell.Value2 = $newValue
}
